$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the last populated row of the data table (column A) and the row
# right below it, where the newest Adafruit IO reading needs to be appended.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

# Copy the previous row down as a starting point so the new cells inherit the
# same (text) cell typing/formatting as the rest of the sheet, then fill in
# the new reading's values.
$ws.Range("A" + $lastRow + ":F" + $lastRow).Copy()
$ws.Range("A" + $newRow + ":F" + $newRow).PasteSpecial(-4163)
$excel.CutCopyMode = 0

$newValues = @("2024-09-25T18:06:40Z", "temperature", "25", "N/A", "N/A", "N/A")
for ($col = 1; $col -le 6; $col++) {
    $cell = $ws.Cells.Item($newRow, $col)
    $target = $newValues[$col - 1]
    if ($cell.Text -ne $target) {
        $cell.Value = $target
    }
}
